$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02801550764843
$ws.Range("D2").Value = 1.030533570490888
$ws.Range("E2").Value = 1.041147292114559
$ws.Range("F2").Value = 1.047901067784118
$ws.Range("I2").Value = 1.030379935605026
$ws.Range("J2").Value = 1.033170172321515
$ws.Range("K2").Value = 1.033344305584337
$ws.Range("L2").Value = 1.043927630590706
$ws.Range("M2").Value = 1.050662411900232
$ws.Range("N2").Value = 1.015064389561419

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028988702716025
$ws.Range("D3").Value = 1.031226782751594
$ws.Range("E3").Value = 1.042141903166911
$ws.Range("F3").Value = 1.049039441223654
$ws.Range("I3").Value = 1.030531646529042
$ws.Range("J3").Value = 1.033783675680557
$ws.Range("K3").Value = 1.033846202007831
$ws.Range("L3").Value = 1.044732326480766
$ws.Range("M3").Value = 1.051611877179309
$ws.Range("N3").Value = 1.015268412893606

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029618448006226
$ws.Range("D4").Value = 1.031674817416383
$ws.Range("E4").Value = 1.042785924746545
$ws.Range("F4").Value = 1.049776784682587
$ws.Range("I4").Value = 1.030627791595851
$ws.Range("J4").Value = 1.034180089201102
$ws.Range("K4").Value = 1.034169769491996
$ws.Range("L4").Value = 1.045252848801738
$ws.Range("M4").Value = 1.05222640139233
$ws.Range("N4").Value = 1.015400202719807

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029883198009453
$ws.Range("D5").Value = 1.031863045528679
$ws.Range("E5").Value = 1.043056776501514
$ws.Range("F5").Value = 1.050086940514234
$ws.Range("I5").Value = 1.030667726259716
$ws.Range("J5").Value = 1.034346605773229
$ws.Range("K5").Value = 1.034305510793863
$ws.Range("L5").Value = 1.045471635153096
$ws.Range("M5").Value = 1.052484784367445
$ws.Range("N5").Value = 1.015455552597609

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029927650985669
$ws.Range("D6").Value = 1.03189464245407
$ws.Range("E6").Value = 1.043102259845678
$ws.Range("F6").Value = 1.05013902738838
$ws.Range("I6").Value = 1.030674403022404
$ws.Range("J6").Value = 1.034374556664805
$ws.Range("K6").Value = 1.034328285532315
$ws.Range("L6").Value = 1.045508367908831
$ws.Range("M6").Value = 1.052528170160341
$ws.Range("N6").Value = 1.015464842878856

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029621985589288
$ws.Range("D7").Value = 1.03167733302488
$ws.Range("E7").Value = 1.04278954346889
$ws.Range("F7").Value = 1.049780928306475
$ws.Range("I7").Value = 1.030628327110229
$ws.Range("J7").Value = 1.034182314737559
$ws.Range("K7").Value = 1.034171584400622
$ws.Range("L7").Value = 1.045255772398448
$ws.Range("M7").Value = 1.052229853773573
$ws.Range("N7").Value = 1.015400942522294

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028344398460358
$ws.Range("D8").Value = 1.030767951795589
$ws.Range("E8").Value = 1.041483333716042
$ws.Range("F8").Value = 1.048285632983888
$ws.Range("I8").Value = 1.030431625317149
$ws.Range("J8").Value = 1.033377625107292
$ws.Range("K8").Value = 1.033514170355451
$ws.Range("L8").Value = 1.044199616118896
$ws.Range("M8").Value = 1.050983255384382
$ws.Range("N8").Value = 1.015133387040635

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026093321070607
$ws.Range("D9").Value = 1.02916157402011
$ws.Range("E9").Value = 1.039185037169734
$ws.Range("F9").Value = 1.045656416255829
$ws.Range("I9").Value = 1.030069553008151
$ws.Range("J9").Value = 1.031955364565717
$ws.Range("K9").Value = 1.032346622192464
$ws.Range("L9").Value = 1.042337259833306
$ws.Range("M9").Value = 1.048787812022154
$ws.Range("N9").Value = 1.014660192738314

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02459275520844
$ws.Range("D10").Value = 1.028088078764859
$ws.Range("E10").Value = 1.037655169923062
$ws.Range("F10").Value = 1.043907450235778
$ws.Range("I10").Value = 1.029817814409285
$ws.Range("J10").Value = 1.031004337078451
$ws.Range("K10").Value = 1.031562185936144
$ws.Range("L10").Value = 1.04109486517524
$ws.Range("M10").Value = 1.047325032302565
$ws.Range("N10").Value = 1.0143435820835

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023943033455178
$ws.Range("D11").Value = 1.027622645834081
$ws.Range("E11").Value = 1.036993280169883
$ws.Range("F11").Value = 1.043151044580941
$ws.Range("I11").Value = 1.029706359588645
$ws.Range("J11").Value = 1.030591861305111
$ws.Range("K11").Value = 1.031221086022722
$ws.Range("L11").Value = 1.040556705371892
$ws.Range("M11").Value = 1.046691838904943
$ws.Range("N11").Value = 1.014206216611726

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023701703041304
$ws.Range("D12").Value = 1.027449673842474
$ws.Range("E12").Value = 1.036747508374488
$ws.Range("F12").Value = 1.042870218496634
$ws.Range("I12").Value = 1.029664592701441
$ws.Range("J12").Value = 1.030438548538854
$ws.Range("K12").Value = 1.031094171608916
$ws.Range("L12").Value = 1.040356780019773
$ws.Range("M12").Value = 1.046456673016479
$ws.Range("N12").Value = 1.014155152441185

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023753468997107
$ws.Range("D13").Value = 1.027486780954517
$ws.Range("E13").Value = 1.036800223471515
$ws.Range("F13").Value = 1.042930450457697
$ws.Range("I13").Value = 1.029673568466911
$ws.Range("J13").Value = 1.030471439206312
$ws.Range("K13").Value = 1.031121404886811
$ws.Range("L13").Value = 1.040399665984897
$ws.Range("M13").Value = 1.046507115520393
$ws.Range("N13").Value = 1.014166107709574

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023923084901233
$ws.Range("D14").Value = 1.027608349734713
$ws.Range("E14").Value = 1.036972962889931
$ws.Range("F14").Value = 1.043127828620269
$ws.Range("I14").Value = 1.029702914620679
$ws.Range("J14").Value = 1.030579190479729
$ws.Range("K14").Value = 1.03121059961386
$ws.Range("L14").Value = 1.040540180065853
$ws.Range("M14").Value = 1.046672399398725
$ws.Range("N14").Value = 1.014201996456912

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024027591549109
$ws.Range("D15").Value = 1.027683240446337
$ws.Range("E15").Value = 1.037079404429178
$ws.Range("F15").Value = 1.043249457934946
$ws.Range("I15").Value = 1.029720947059245
$ws.Range("J15").Value = 1.030645566238925
$ws.Range("K15").Value = 1.031265526997502
$ws.Range("L15").Value = 1.040626751619496
$ws.Range("M15").Value = 1.046774240279591
$ws.Range("N15").Value = 1.014224103335293

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024635875816555
$ws.Range("D16").Value = 1.028118955411524
$ws.Range("E16").Value = 1.037699109054303
$ws.Range("F16").Value = 1.043957669540237
$ws.Range("I16").Value = 1.029825159718284
$ws.Range("J16").Value = 1.031031697562782
$ws.Range("K16").Value = 1.031584793464573
$ws.Range("L16").Value = 1.041130577001987
$ws.Range("M16").Value = 1.047367059500075
$ws.Range("N16").Value = 1.014352692884872

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025017445482477
$ws.Range("D17").Value = 1.028392107345426
$ws.Range("E17").Value = 1.038087981741685
$ws.Range("F17").Value = 1.044402155263805
$ws.Range("I17").Value = 1.029889873916782
$ws.Range("J17").Value = 1.031273727252854
$ws.Range("K17").Value = 1.031784677501625
$ws.Range("L17").Value = 1.04144656151598
$ws.Range("M17").Value = 1.047738973062213
$ws.Range("N17").Value = 1.014433281301981

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025240011800449
$ws.Range("D18").Value = 1.028551374055836
$ws.Range("E18").Value = 1.038314857912812
$ws.Range("F18").Value = 1.044661503877708
$ws.Range("I18").Value = 1.029927384155663
$ws.Range("J18").Value = 1.031414833957522
$ws.Range("K18").Value = 1.03190112808762
$ws.Range("L18").Value = 1.04163085112192
$ws.Range("M18").Value = 1.047955923215458
$ws.Range("N18").Value = 1.014480261004886

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025315901698051
$ws.Range("D19").Value = 1.028605669986651
$ws.Range("E19").Value = 1.038392225852454
$ws.Range("F19").Value = 1.044749949906819
$ws.Range("I19").Value = 1.029940134054059
$ws.Range("J19").Value = 1.031462936638214
$ws.Range("K19").Value = 1.031940811233177
$ws.Range("L19").Value = 1.041693685911052
$ws.Range("M19").Value = 1.048029900869578
$ws.Range("N19").Value = 1.01449627542691

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024976506337267
$ws.Range("D20").Value = 1.028362806739503
$ws.Range("E20").Value = 1.038046253864806
$ws.Range("F20").Value = 1.044354457088515
$ws.Range("I20").Value = 1.029882955149817
$ws.Range("J20").Value = 1.031247766489846
$ws.Range("K20").Value = 1.031763246145435
$ws.Range("L20").Value = 1.041412661311112
$ws.Range("M20").Value = 1.047699068261805
$ws.Range("N20").Value = 1.014424637631671

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023873137100676
$ws.Range("D21").Value = 1.027572553220795
$ws.Range("E21").Value = 1.036922093127052
$ws.Range("F21").Value = 1.043069701896874
$ws.Range("I21").Value = 1.02969428305326
$ws.Range("J21").Value = 1.030547463197024
$ws.Range("K21").Value = 1.031184339911565
$ws.Range("L21").Value = 1.040498802964034
$ws.Range("M21").Value = 1.046623726576401
$ws.Range("N21").Value = 1.01419142923229

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023179435098394
$ws.Range("D22").Value = 1.027075172279327
$ws.Range("E22").Value = 1.036215772159452
$ws.Range("F22").Value = 1.042262716062901
$ws.Range("I22").Value = 1.029573530831871
$ws.Range("J22").Value = 1.030106571174589
$ws.Range("K22").Value = 1.030819116631527
$ws.Range("L22").Value = 1.039924057296119
$ws.Range("M22").Value = 1.0459477922208
$ws.Range("N22").Value = 1.014044567524054

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023547176745372
$ws.Range("D23").Value = 1.027338892010531
$ws.Range("E23").Value = 1.036590160282067
$ws.Range("F23").Value = 1.042690439498169
$ws.Range("I23").Value = 1.029637745268745
$ws.Range("J23").Value = 1.030340351541463
$ws.Range("K23").Value = 1.031012845893873
$ws.Range("L23").Value = 1.040228756508636
$ws.Range("M23").Value = 1.046306101062284
$ws.Range("N23").Value = 1.014122443842585

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024995004964824
$ws.Range("D24").Value = 1.028376046601937
$ws.Range("E24").Value = 1.038065108731429
$ws.Range("F24").Value = 1.044376009569964
$ws.Range("I24").Value = 1.029886082173422
$ws.Range("J24").Value = 1.031259497242353
$ws.Range("K24").Value = 1.031772930480404
$ws.Range("L24").Value = 1.041427979412553
$ws.Range("M24").Value = 1.047717099463813
$ws.Range("N24").Value = 1.014428543415039

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026675253251066
$ws.Range("D25").Value = 1.029577320201899
$ws.Range("E25").Value = 1.039778793812731
$ws.Range("F25").Value = 1.046335455750526
$ws.Range("I25").Value = 1.030164985523665
$ws.Range("J25").Value = 1.032323558452468
$ws.Range("K25").Value = 1.032649534878781
$ws.Range("L25").Value = 1.042818871155482
$ws.Range("M25").Value = 1.04935523907354
$ws.Range("N25").Value = 1.014782728101261
